$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 199.71428
$ws.Range("I11").Value = 199.71428
$ws.Range("K11").Value = 199.71428
$ws.Range("M11").Value = -59.71428

$ws.Range("H33").Value = 305
$ws.Range("I33").Value = 94.70587999999999
$ws.Range("J33").Value = 815.7143
$ws.Range("K33").Value = 94.70587999999999
$ws.Range("L33").Value = 815.7143
$ws.Range("M33").Value = 134.29412
$ws.Range("N33").Value = -1273.7143

$ws.Range("H76").Value = 4871.826
$ws.Range("I76").Value = 3500.125
$ws.Range("J76").Value = 5603.4
$ws.Range("K76").Value = 3500.125
$ws.Range("L76").Value = 5603.4
$ws.Range("M76").Value = -3185.125
$ws.Range("N76").Value = -6233.4

$ws.Range("H79").Value = 4871.826
$ws.Range("I79").Value = 3500.125
$ws.Range("J79").Value = 5603.4
$ws.Range("K79").Value = 3500.125
$ws.Range("L79").Value = 5603.4
$ws.Range("M79").Value = -2408.125
$ws.Range("N79").Value = -7787.4

$ws.Range("H80").Value = 54708.785
$ws.Range("I80").Value = 91699.73
$ws.Range("J80").Value = 455.4
$ws.Range("K80").Value = 275099.19
$ws.Range("L80").Value = 1366.2
$ws.Range("M80").Value = -274101.19
$ws.Range("N80").Value = -3362.2

$ws.Range("H83").Value = 54708.785
$ws.Range("I83").Value = 91699.73
$ws.Range("J83").Value = 455.4
$ws.Range("K83").Value = 825297.5699999999
$ws.Range("L83").Value = 4098.599999999999
$ws.Range("M83").Value = -820305.5699999999
$ws.Range("N83").Value = -14082.6

$ws.Range("H115").Value = 984.7778
$ws.Range("I115").Value = 623.2857
$ws.Range("J115").Value = 2250
$ws.Range("K115").Value = 1869.8571
$ws.Range("L115").Value = 6750
$ws.Range("M115").Value = -302.8571000000002
$ws.Range("N115").Value = -9884

$ws.Range("H132").Value = 3510.8865
$ws.Range("I132").Value = 3726.8
$ws.Range("J132").Value = 2671.2222
$ws.Range("K132").Value = 11180.4
$ws.Range("L132").Value = 8013.6666
$ws.Range("M132").Value = -8650.400000000001
$ws.Range("N132").Value = -13073.6666

$ws.Range("H137").Value = 10405.4
$ws.Range("I137").Value = 22571.4
$ws.Range("J137").Value = 2294.7334
$ws.Range("K137").Value = 67714.20000000001
$ws.Range("L137").Value = 6884.2002
$ws.Range("M137").Value = -65164.20000000001
$ws.Range("N137").Value = -11984.2002

$ws.Range("H138").Value = 1475.2778
$ws.Range("I138").Value = 838.84
$ws.Range("J138").Value = 2921.7273
$ws.Range("K138").Value = 2516.52
$ws.Range("L138").Value = 8765.1819
$ws.Range("M138").Value = 2623.48
$ws.Range("N138").Value = -19045.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 12000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 12000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 12000
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -12700

$ws.Range("H45").Value = 53956.074
$ws.Range("I45").Value = 109187.16
$ws.Range("J45").Value = 3985.0952
$ws.Range("K45").Value = 109187.16
$ws.Range("L45").Value = 3985.0952
$ws.Range("M45").Value = -108810.16
$ws.Range("N45").Value = -4739.0952

$ws.Range("H97").Value = 13340797
$ws.Range("I97").Value = 13220.125
$ws.Range("J97").Value = 28572314
$ws.Range("K97").Value = 13220.125
$ws.Range("L97").Value = 28572314
$ws.Range("M97").Value = -12724.125
$ws.Range("N97").Value = -28573306

$ws.Range("H102").Value = 6782.816
$ws.Range("I102").Value = 8032.3667
$ws.Range("J102").Value = 2097
$ws.Range("K102").Value = 8032.3667
$ws.Range("L102").Value = 2097
$ws.Range("M102").Value = -6410.3667
$ws.Range("N102").Value = -5341

$ws.Range("H122").Value = 2150177.5
$ws.Range("J122").Value = 3338780.5
$ws.Range("L122").Value = 10016341.5
$ws.Range("N122").Value = -10021241.5

$ws.Range("H132").Value = 3727.24
$ws.Range("I132").Value = 4094.875
$ws.Range("J132").Value = 3073.6667
$ws.Range("K132").Value = 12284.625
$ws.Range("L132").Value = 9221.000100000001
$ws.Range("M132").Value = -9754.625
$ws.Range("N132").Value = -14281.0001

$ws.Range("H139").Value = 233300
$ws.Range("J139").Value = 233300
$ws.Range("L139").Value = 233300
$ws.Range("N139").Value = -243580

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 15000
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

$ws.Range("H82").Value = 12072.857
$ws.Range("J82").Value = 60000
$ws.Range("L82").Value = 60000
$ws.Range("N82").Value = -60766

$ws.Range("H85").Value = 12072.857
$ws.Range("J85").Value = 60000
$ws.Range("L85").Value = 60000
$ws.Range("N85").Value = -62652

$ws.Range("H94").Value = 9145.677
$ws.Range("I94").Value = 12699.5
$ws.Range("J94").Value = 2630.3333
$ws.Range("K94").Value = 12699.5
$ws.Range("L94").Value = 2630.3333
$ws.Range("M94").Value = -12248.5
$ws.Range("N94").Value = -3532.3333

$ws.Range("H134").Value = 7516.909
$ws.Range("I134").Value = 8715.166999999999
$ws.Range("J134").Value = 2124.75
$ws.Range("K134").Value = 26145.501
$ws.Range("L134").Value = 6374.25
$ws.Range("M134").Value = -23610.501
$ws.Range("N134").Value = -11444.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19619.875
$ws.Range("J31").Value = 5662.1665
$ws.Range("L31").Value = 5662.1665
$ws.Range("N31").Value = -6252.1665

$ws.Range("H34").Value = 19619.875
$ws.Range("J34").Value = 5662.1665
$ws.Range("L34").Value = 5662.1665
$ws.Range("N34").Value = -6066.1665

$ws.Range("H132").Value = 3999.7144
$ws.Range("I132").Value = 3999.6
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 11998.8
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -9468.799999999999
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 9.25
$ws.Range("I26").Value = 1
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 285

$ws.Range("H46").Value = 4240.9
$ws.Range("I46").Value = 887.1667
$ws.Range("J46").Value = 5678.2144
$ws.Range("K46").Value = 2661.5001
$ws.Range("L46").Value = 17034.6432
$ws.Range("M46").Value = -2570.5001
$ws.Range("N46").Value = -17216.6432

$ws.Range("H113").Value = 26643.75
$ws.Range("I113").Value = 837.5
$ws.Range("J113").Value = 52450
$ws.Range("K113").Value = 2512.5
$ws.Range("L113").Value = 157350
$ws.Range("M113").Value = -342.5
$ws.Range("N113").Value = -161690

$ws.Range("H129").Value = 2522.6667
$ws.Range("J129").Value = 3171.818
$ws.Range("L129").Value = 9515.454000000002
$ws.Range("N129").Value = -19515.454

$ws.Range("H136").Value = 2379.3125
$ws.Range("I136").Value = 1291.0714
$ws.Range("K136").Value = 3873.2142
$ws.Range("M136").Value = 1226.7858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7801.278
$ws.Range("I70").Value = 6358.231
$ws.Range("K70").Value = 6358.231
$ws.Range("M70").Value = -6088.231

$ws.Range("H73").Value = 7801.278
$ws.Range("I73").Value = 6358.231
$ws.Range("K73").Value = 6358.231
$ws.Range("M73").Value = -5422.231

$ws.Range("H102").Value = 10418.333
$ws.Range("I102").Value = 14096.272
$ws.Range("K102").Value = 14096.272
$ws.Range("M102").Value = -12474.272

$ws.Range("H132").Value = 4031.5435
$ws.Range("I132").Value = 4498.647
$ws.Range("K132").Value = 13495.941
$ws.Range("M132").Value = -10965.941

$ws.Range("H136").Value = 33644.832
$ws.Range("J136").Value = 33644.832
$ws.Range("L136").Value = 100934.496
$ws.Range("N136").Value = -106034.496

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3692.3635
$ws.Range("I82").Value = 3852
$ws.Range("J82").Value = 3266.6667
$ws.Range("K82").Value = 3852
$ws.Range("L82").Value = 3266.6667
$ws.Range("M82").Value = -3491
$ws.Range("N82").Value = -3988.6667

$ws.Range("H85").Value = 3692.3635
$ws.Range("I85").Value = 3852
$ws.Range("J85").Value = 3266.6667
$ws.Range("K85").Value = 3852
$ws.Range("L85").Value = 3266.6667
$ws.Range("M85").Value = -2604
$ws.Range("N85").Value = -5762.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 1259500
$ws.Range("I24").Value = 1259500
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 1259500
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -1259270
$ws.Range("N24").ClearContents()

$ws.Range("H30").Value = 83333.336
$ws.Range("J30").Value = 50000
$ws.Range("L30").Value = 50000
$ws.Range("N30").Value = -50214

$ws.Range("H95").Value = 172020180
$ws.Range("J95").Value = 172020180
$ws.Range("L95").Value = 172020180
$ws.Range("N95").Value = -172025672

$ws.Range("H122").Value = 4265.0215
$ws.Range("I122").Value = 2133.879
$ws.Range("K122").Value = 6401.637
$ws.Range("M122").Value = -3951.637

$ws.Range("H132").Value = 11966.625
$ws.Range("I132").Value = 15863.704
$ws.Range("J132").Value = 3872.6924
$ws.Range("K132").Value = 47591.112
$ws.Range("L132").Value = 11618.0772
$ws.Range("M132").Value = -45061.112
$ws.Range("N132").Value = -16678.0772

$ws.Range("H137").Value = 62333.332
$ws.Range("J137").Value = 62333.332
$ws.Range("L137").Value = 62333.332
$ws.Range("N137").Value = -72533.33199999999

Write-Output "Applied all changes"